{"js": "// The Jekyll site rebuild that produced this change strips the page's\n// trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" action line and\n// the \"\u00a9 2020 ... Jekyll and Github pages ...\" footer line (plus the blank\n// spacer paragraph that separated them from the body), right after the\n// \"LOM3212: Fen\u00f4menos de Transporte A (Requisito)\" requirement line.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the requirement paragraph that anchors the block being removed.\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"LOM3212\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOM3212' anchor paragraph\");\n}\n\n// Delete the 3 paragraphs that immediately follow the anchor: the blank\n// spacer, the \"Ver no Jupiter...\" line, and the site footer/copyright line.\n// Delete from the highest index down so earlier deletions don't invalidate\n// the indices of paragraphs still waiting to be removed.\nfor (let offset = 3; offset >= 1; offset--) {\n  paragraphs.items[anchorIndex + offset].delete();\n}\n\nawait context.sync();\n", "ps1": "# The Jekyll site rebuild that produced this change strips the page's\n# trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" action line and\n# the \"(c) 2020 ... Jekyll and Github pages ...\" footer line (plus the blank\n# spacer paragraph that separated them from the body), right after the\n# \"LOM3212: Fenomenos de Transporte A (Requisito)\" requirement line.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"LOM3212\")\nif (-not $found) {\n    throw \"Could not find the 'LOM3212' anchor paragraph\"\n}\n$anchorStart = $rng.Start\n$anchorEnd = $rng.End\n\n# Resolve which paragraph in the document's Paragraphs collection the match\n# landed in (Find only returns the matched substring, not the whole\n# paragraph), so we can index relative to it.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $anchorStart -and $p.Range.End -ge $anchorEnd) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve the anchor paragraph index\"\n}\n\n# Delete the 3 paragraphs that immediately follow the anchor: the blank\n# spacer, the \"Ver no Jupiter...\" line, and the site footer/copyright line.\n# Delete from the highest index down so earlier deletions don't shift the\n# indices of paragraphs still waiting to be removed.\n$d.Paragraphs.Item($anchorIndex + 3).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 2).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n"}
